$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71; existing rows 71-142 shift down to 72-143.
$ws.Rows(71).Insert()

# Populate the newly inserted row 71 with the new data record.
$ws.Range("A71").Value = 5
$ws.Range("B71").Value = "Macroferia Regional de Talca"
$ws.Range("C71").Value = "Maule"
$ws.Range("D71").Value = 44539
$ws.Range("E71").Value = 7
$ws.Range("F71").Value = 100112021
$ws.Range("G71").Value = "Ají"
$ws.Range("H71").Value = "Americana (o)"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 60
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = 15000
$ws.Range("N71").Value = "`$/caja 15 kilos"
$ws.Range("O71").Value = "Región del Maule"
$ws.Range("P71").Value = 1000
$ws.Range("Q71").Value = 15
$ws.Range("R71").Value = "Hortaliza"
